$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells for the new columns, matching the style of the existing header row (col AC, style index 1)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header style (bold font, border, centered alignment) from an existing header cell
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the Wins / Losses / Ties data for each data row (2-50)
for ($r = 2; $r -le 50; $r++) {
    $ws.Cells.Item($r, 30).Value = 90   # AD
    $ws.Cells.Item($r, 31).Value = 72   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
